$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1 (clause SEXTA): remove " hasta $fecha_renuncia" so the sentence
# ends at "...finalización de la relación." — the trailing period stays as
# its own run, and the stray proofErr markers that wrapped "fecha_renuncia"
# are removed along with it.
# ---------------------------------------------------------------------------
$find1 = $d.Content.Duplicate
$found1 = $find1.Find.Execute(" hasta " + [char]36 + "fecha_renuncia.")

if ($found1) {
    $target1 = $d.Range($find1.Start, $find1.End)
    $innerXml1 = "<w:r $w><w:t>.</w:t></w:r>"
    $xmlFrag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part1.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml1 + '</w:p></pkg:xmlData></pkg:part></pkg:package>'
    $target1.InsertXML($xmlFrag1)
}

# ---------------------------------------------------------------------------
# Change 2 (Sanciones paragraph): replace the "$cantidad_sancion" merge
# field placeholder with a literal amount: "$100.00 Dólares de los Estados
# Unidos de América".
# ---------------------------------------------------------------------------
$dolares = "D" + [char]0x00F3 + "lares"
$america = " de los Estados Unidos de Am" + [char]0x00E9 + "rica"

$find2 = $d.Content.Duplicate
$found2 = $find2.Find.Execute([char]36 + "cantidad_sancion.")

if ($found2) {
    $target2 = $d.Range($find2.Start, $find2.End)
    $innerXml2 = "<w:r $w><w:t>" + [char]36 + "</w:t></w:r>" +
                 "<w:r $w><w:t xml:space=`"preserve`">100.00 </w:t></w:r>" +
                 "<w:proofErr $w w:type=`"gramStart`"/>" +
                 "<w:r $w><w:t>$dolares</w:t></w:r>" +
                 "<w:proofErr $w w:type=`"gramEnd`"/>" +
                 "<w:r $w><w:t xml:space=`"preserve`">$america</w:t></w:r>" +
                 "<w:r $w><w:t>.</w:t></w:r>"
    $xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part2.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml2 + '</w:p></pkg:xmlData></pkg:part></pkg:package>'
    $target2.InsertXML($xmlFrag2)
}
